$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1623.0555
$ws.Range("J39").Value = 1191.2858
$ws.Range("L39").Value = 3573.8574
$ws.Range("N39").Value = -4165.857400000001
$ws.Range("H42").Value = 1039.3334
$ws.Range("I42").Value = 100
$ws.Range("K42").Value = 300
$ws.Range("M42").Value = -70
$ws.Range("H58").Value = 1013.5714
$ws.Range("I58").Value = 273.75
$ws.Range("K58").Value = 821.25
$ws.Range("M58").Value = -671.25
$ws.Range("H107").Value = 40224.42
$ws.Range("I107").Value = 541.5
$ws.Range("K107").Value = 541.5
$ws.Range("M107").Value = 1378.5
$ws.Range("H137").Value = 2222.7097
$ws.Range("I137").Value = 1295.0555
$ws.Range("J137").Value = 3507.1538
$ws.Range("K137").Value = 3885.1665
$ws.Range("L137").Value = 10521.4614
$ws.Range("M137").Value = -1335.1665
$ws.Range("N137").Value = -15621.4614
$ws.Range("H138").Value = 2879.6453
$ws.Range("J138").Value = 3360.6
$ws.Range("L138").Value = 10081.8
$ws.Range("N138").Value = -20361.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3167.3333
$ws.Range("I2").Value = 2359.5715
$ws.Range("J2").Value = 5994.5
$ws.Range("K2").Value = 2359.5715
$ws.Range("L2").Value = 5994.5
$ws.Range("M2").Value = -2246.5715
$ws.Range("N2").Value = -6220.5
$ws.Range("H45").Value = 8706.777
$ws.Range("I45").Value = 14687.125
$ws.Range("K45").Value = 14687.125
$ws.Range("M45").Value = -14310.125
$ws.Range("H57").Value = 5171.4287
$ws.Range("I57").Value = 5171.4287
$ws.Range("K57").Value = 5171.4287
$ws.Range("M57").Value = -4687.4287
$ws.Range("H61").Value = 2036.7192
$ws.Range("I61").Value = 1978.9524
$ws.Range("K61").Value = 1978.9524
$ws.Range("M61").Value = -1766.9524
$ws.Range("H74").Value = 2086.32
$ws.Range("I74").Value = 1811.5
$ws.Range("J74").Value = 2793
$ws.Range("K74").Value = 1811.5
$ws.Range("L74").Value = 2793
$ws.Range("M74").Value = -937.5
$ws.Range("N74").Value = -4541
$ws.Range("H77").Value = 2086.32
$ws.Range("I77").Value = 1811.5
$ws.Range("J77").Value = 2793
$ws.Range("K77").Value = 9057.5
$ws.Range("L77").Value = 13965
$ws.Range("M77").Value = -4689.5
$ws.Range("N77").Value = -22701
$ws.Range("H116").Value = 3167.3333
$ws.Range("I116").Value = 2359.5715
$ws.Range("J116").Value = 5994.5
$ws.Range("K116").Value = 2359.5715
$ws.Range("L116").Value = 5994.5
$ws.Range("M116").Value = -65.57150000000001
$ws.Range("N116").Value = -10582.5
$ws.Range("H126").Value = 5068.3076
$ws.Range("I126").Value = 5068.3076
$ws.Range("K126").Value = 15204.9228
$ws.Range("M126").Value = -12734.9228
$ws.Range("H132").Value = 1477.4043
$ws.Range("I132").Value = 1476.4222
$ws.Range("K132").Value = 4429.2666
$ws.Range("M132").Value = -1899.2666
$ws.Range("H136").Value = 2036.7192
$ws.Range("I136").Value = 1978.9524
$ws.Range("K136").Value = 5936.857199999999
$ws.Range("M136").Value = -3386.857199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3167.3333
$ws.Range("I3").Value = 2359.5715
$ws.Range("J3").Value = 5994.5
$ws.Range("K3").Value = 2359.5715
$ws.Range("L3").Value = 5994.5
$ws.Range("M3").Value = -2245.5715
$ws.Range("N3").Value = -6222.5
$ws.Range("H22").Value = 758
$ws.Range("I22").Value = 758
$ws.Range("K22").Value = 758
$ws.Range("M22").Value = -585
$ws.Range("H107").Value = 49409.715
$ws.Range("I107").Value = 57383.555
$ws.Range("J107").Value = 1566.6666
$ws.Range("K107").Value = 57383.555
$ws.Range("L107").Value = 1566.6666
$ws.Range("M107").Value = -55463.555
$ws.Range("N107").Value = -5406.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1953.6666
$ws.Range("I16").Value = 1669.7142
$ws.Range("K16").Value = 1669.7142
$ws.Range("M16").Value = -1382.7142
$ws.Range("H23").Value = 8000
$ws.Range("J23").Value = 8000
$ws.Range("L23").Value = 8000
$ws.Range("N23").Value = -8480
$ws.Range("H27").Value = 8000
$ws.Range("J27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("N27").Value = -8384
$ws.Range("H58").Value = 2349.75
$ws.Range("I58").Value = 2145.25
$ws.Range("K58").Value = 2145.25
$ws.Range("M58").Value = -1942.25
$ws.Range("H107").Value = 1990.2069
$ws.Range("I107").Value = 1813
$ws.Range("K107").Value = 1813
$ws.Range("M107").Value = 107
$ws.Range("H113").Value = 1953.6666
$ws.Range("I113").Value = 1669.7142
$ws.Range("K113").Value = 1669.7142
$ws.Range("M113").Value = 500.2858000000001
$ws.Range("H136").Value = 2349.75
$ws.Range("I136").Value = 2145.25
$ws.Range("K136").Value = 6435.75
$ws.Range("M136").Value = -3885.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 703.2414
$ws.Range("J107").Value = 689
$ws.Range("L107").Value = 2067
$ws.Range("N107").Value = -5907
$ws.Range("H116").Value = 5000
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 15000
$ws.Range("N116").Value = -21884
$ws.Range("H117").Value = 5330.278
$ws.Range("I117").Value = 1342.375
$ws.Range("J117").Value = 8520.6
$ws.Range("K117").Value = 4027.125
$ws.Range("L117").Value = 25561.8
$ws.Range("M117").Value = -585.125
$ws.Range("N117").Value = -32445.8
$ws.Range("H119").Value = 1015.25
$ws.Range("I119").Value = 1262.3334
$ws.Range("J119").Value = 274
$ws.Range("K119").Value = 3787.0002
$ws.Range("L119").Value = 822
$ws.Range("M119").Value = 1050.9998
$ws.Range("N119").Value = -10498
$ws.Range("H120").Value = 11500
$ws.Range("J120").Value = 20000
$ws.Range("L120").Value = 60000
$ws.Range("N120").Value = -69676
$ws.Range("H137").Value = 4027.4666
$ws.Range("J137").Value = 4635.636
$ws.Range("L137").Value = 13906.908
$ws.Range("N137").Value = -24106.908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2333871.5
$ws.Range("I3").Value = 2200069.5
$ws.Range("J3").Value = 2501123.8
$ws.Range("K3").Value = 2200069.5
$ws.Range("L3").Value = 2501123.8
$ws.Range("M3").Value = -2199953.5
$ws.Range("N3").Value = -2501355.8
$ws.Range("H41").Value = 7663.3335
$ws.Range("I41").Value = 6495
$ws.Range("K41").Value = 6495
$ws.Range("M41").Value = -6140
$ws.Range("H80").Value = 3370.2307
$ws.Range("I80").Value = 2507.1428
$ws.Range("J80").Value = 4377.1665
$ws.Range("K80").Value = 2507.1428
$ws.Range("L80").Value = 4377.1665
$ws.Range("M80").Value = -1509.1428
$ws.Range("N80").Value = -6373.1665
$ws.Range("H83").Value = 3370.2307
$ws.Range("I83").Value = 2507.1428
$ws.Range("J83").Value = 4377.1665
$ws.Range("K83").Value = 12535.714
$ws.Range("L83").Value = 21885.8325
$ws.Range("M83").Value = -7543.714
$ws.Range("N83").Value = -31869.8325
$ws.Range("H132").Value = 6538.4707
$ws.Range("I132").Value = 6263.9
$ws.Range("K132").Value = 18791.7
$ws.Range("M132").Value = -16261.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 477.42105
$ws.Range("I55").Value = 557.75
$ws.Range("K55").Value = 557.75
$ws.Range("M55").Value = -384.75
$ws.Range("H68").Value = 4029.8572
$ws.Range("I68").Value = 3674.4546
$ws.Range("J68").Value = 5333
$ws.Range("K68").Value = 3674.4546
$ws.Range("L68").Value = 5333
$ws.Range("M68").Value = -2925.4546
$ws.Range("N68").Value = -6831
$ws.Range("H71").Value = 4029.8572
$ws.Range("I71").Value = 3674.4546
$ws.Range("J71").Value = 5333
$ws.Range("K71").Value = 18372.273
$ws.Range("L71").Value = 26665
$ws.Range("M71").Value = -14628.273
$ws.Range("N71").Value = -34153
$ws.Range("H122").Value = 94951.87
$ws.Range("I122").Value = 128286
$ws.Range("J122").Value = 6060.8335
$ws.Range("K122").Value = 384858
$ws.Range("L122").Value = 18182.5005
$ws.Range("M122").Value = -382408
$ws.Range("N122").Value = -23082.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2726219.5
$ws.Range("I132").Value = 1720457.9
$ws.Range("J132").Value = 5296499
$ws.Range("K132").Value = 5161373.699999999
$ws.Range("L132").Value = 15889497
$ws.Range("M132").Value = -5158843.699999999
$ws.Range("N132").Value = -15894557
